$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the current column F ("GFA - Sales"),
# shifting the old F:M columns to H:O. Excel carries formatting/values
# of the existing column F into the new slots, so insert at F twice.
$ws.Range("F:G").Insert()

# Header for the two new columns (style is already inherited from the insert)
$ws.Cells.Item(1, 6).Value = "M_TotalTax"
$ws.Cells.Item(1, 7).Value = "M_CorpTax"

# New data values for the two inserted columns (rows 2-11)
$newValues = @{
    2  = @(6308727034979.312, 399825921028.5854)
    3  = @(16630145391623.02, 1639742485782.957)
    4  = @(4450994137606.095, 601350231413.5104)
    5  = @(4183547438952.192, 598849276038.3025)
    6  = @(11223287075501.79, 872292028558.4308)
    7  = @(1841737275230.086, 214321200777.9413)
    8  = @(6192585801479.285, 516695167857.3162)
    9  = @(14653861967257.56, 1232540278767.842)
    10 = @(9623160693235.053, 876943418066.7275)
    11 = @(5030701274022.499, 355596860701.1148)
}

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    $ws.Cells.Item($row, 6).Value = $vals[0]
    $ws.Cells.Item($row, 7).Value = $vals[1]
}
